$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column P for rows 2-14 (data removed in Sept 24 SSP update)
$ws.Range("P2:P14").ClearContents()

# Update columns Q:AA for rows 2-13 with new values
$ws.Range("Q2").Value = 34.69638015858982
$ws.Range("R2").Value = -26.9389830298353
$ws.Range("S2").Value = -146.4574609384943
$ws.Range("T2").Value = -283.8968469856518
$ws.Range("U2").Value = -391.1227627719198
$ws.Range("V2").Value = -490.519535797212
$ws.Range("W2").Value = -812.332694417669
$ws.Range("X2").Value = -723.1289987016471
$ws.Range("Y2").Value = -1996.256981084404
$ws.Range("Z2").Value = -1910.469919733242
$ws.Range("AA2").Value = -1201.69277056063
$ws.Range("Q3").Value = -171.8816458255606
$ws.Range("R3").Value = -766.8029417826619
$ws.Range("S3").Value = -985.1894092423827
$ws.Range("T3").Value = -1057.891732932264
$ws.Range("U3").Value = -996.3330636358082
$ws.Range("V3").Value = -1186.972913996136
$ws.Range("W3").Value = -2465.199073890797
$ws.Range("X3").Value = -3291.656942142021
$ws.Range("Y3").Value = -2851.813327760527
$ws.Range("Z3").Value = -1963.862416315807
$ws.Range("AA3").Value = -1691.236984979456
$ws.Range("Q4").Value = 7.431200764599907
$ws.Range("R4").Value = -22.02977117085794
$ws.Range("S4").Value = -83.11760065520056
$ws.Range("T4").Value = -147.6054858712574
$ws.Range("U4").Value = -220.7351165972183
$ws.Range("V4").Value = -303.7603590478333
$ws.Range("W4").Value = -657.7990617456937
$ws.Range("X4").Value = -1091.420065494662
$ws.Range("Y4").Value = -877.1508935295822
$ws.Range("Z4").Value = -811.5433269917892
$ws.Range("AA4").Value = -351.0280795903536
$ws.Range("Q5").Value = 45.88772144876566
$ws.Range("R5").Value = -0.5035406081292422
$ws.Range("S5").Value = -42.08509019622556
$ws.Range("T5").Value = -80.04191050949044
$ws.Range("U5").Value = -153.1955808062492
$ws.Range("V5").Value = -243.903702765754
$ws.Range("W5").Value = -623.9250360012662
$ws.Range("X5").Value = -1482.857243862202
$ws.Range("Y5").Value = -1729.78595086159
$ws.Range("Z5").Value = -2647.891347494049
$ws.Range("AA5").Value = -5759.994447360386
$ws.Range("Q6").Value = 63.18027215502912
$ws.Range("R6").Value = -7.070858447161299
$ws.Range("S6").Value = -73.50938246215605
$ws.Range("T6").Value = -231.3563518091696
$ws.Range("U6").Value = -354.1100394266708
$ws.Range("V6").Value = -420.3176936272977
$ws.Range("W6").Value = -889.3360357633879
$ws.Range("X6").Value = -2733.081710337864
$ws.Range("Y6").Value = -3632.92547371078
$ws.Range("Z6").Value = -6787.655561917422
$ws.Range("AA6").Value = -13135.84495522096
$ws.Range("Q7").Value = 109.3333343292943
$ws.Range("R7").Value = -55.70858328534996
$ws.Range("S7").Value = -160.4280375032772
$ws.Range("T7").Value = -254.4253579441036
$ws.Range("U7").Value = -340.1986459188897
$ws.Range("V7").Value = -469.664857884751
$ws.Range("W7").Value = -851.0965359428355
$ws.Range("X7").Value = -1073.533883216317
$ws.Range("Y7").Value = -646.3702680292732
$ws.Range("Z7").Value = -289.2864647605688
$ws.Range("AA7").Value = 85.19946350476607
$ws.Range("Q8").Value = 123.7520147028954
$ws.Range("R8").Value = 18.55221283122182
$ws.Range("S8").Value = -68.27148455082533
$ws.Range("T8").Value = -177.5502885604708
$ws.Range("U8").Value = -309.8990866258553
$ws.Range("V8").Value = -543.8861411260295
$ws.Range("W8").Value = -1944.171021147397
$ws.Range("X8").Value = -6798.230133646563
$ws.Range("Y8").Value = -6723.593926639877
$ws.Range("Z8").Value = -6060.89959787308
$ws.Range("AA8").Value = -1220.18463416079
$ws.Range("Q9").Value = 38.17585913825769
$ws.Range("R9").Value = -8.183226591304901
$ws.Range("S9").Value = -48.8635156627721
$ws.Range("T9").Value = -105.6351269118392
$ws.Range("U9").Value = -174.7068081505758
$ws.Range("V9").Value = -236.5809447984768
$ws.Range("W9").Value = -722.7311709081712
$ws.Range("X9").Value = -1604.281185181418
$ws.Range("Y9").Value = -1455.67777630344
$ws.Range("Z9").Value = -1246.170272489308
$ws.Range("AA9").Value = -354.1848721814411
$ws.Range("Q10").Value = -20.99015300688166
$ws.Range("R10").Value = -112.3516688555226
$ws.Range("S10").Value = -200.2629971996486
$ws.Range("T10").Value = -299.1664186074325
$ws.Range("U10").Value = -405.0147993466621
$ws.Range("V10").Value = -513.6682477290864
$ws.Range("W10").Value = -1023.483919056334
$ws.Range("X10").Value = -1528.595673540802
$ws.Range("Y10").Value = -1169.824334542595
$ws.Range("Z10").Value = -995.1323565447985
$ws.Range("AA10").Value = -790.8226253271453
$ws.Range("Q11").Value = 16.00300714702496
$ws.Range("R11").Value = 9.960796898525302
$ws.Range("S11").Value = -34.05020768249531
$ws.Range("T11").Value = -90.99952318390088
$ws.Range("U11").Value = -167.6304709175562
$ws.Range("V11").Value = -344.6622603192554
$ws.Range("W11").Value = -1570.373584514152
$ws.Range("X11").Value = -4925.14461820703
$ws.Range("Y11").Value = -4558.038925791524
$ws.Range("Z11").Value = -3394.863172406475
$ws.Range("AA11").Value = -2088.653241936828
$ws.Range("Q12").Value = -77.89775068598692
$ws.Range("R12").Value = -236.5379925745895
$ws.Range("S12").Value = -457.4003502513031
$ws.Range("T12").Value = -665.2581465679891
$ws.Range("U12").Value = -749.7447725900686
$ws.Range("V12").Value = -815.8900747733443
$ws.Range("W12").Value = -1145.674765810396
$ws.Range("X12").Value = -685.4746981795565
$ws.Range("Y12").Value = -698.0676880510742
$ws.Range("Z12").Value = -1005.934822344592
$ws.Range("AA12").Value = -674.5117732418621
$ws.Range("Q13").Value = 28.91629080955727
$ws.Range("R13").Value = -35.18547274014122
$ws.Range("S13").Value = -143.5563513469347
$ws.Range("T13").Value = -284.5316668555723
$ws.Range("U13").Value = -427.1915298542483
$ws.Range("V13").Value = -668.068460830879
$ws.Range("W13").Value = -1829.924452021452
$ws.Range("X13").Value = -4647.406677756597
$ws.Range("Y13").Value = -4184.342536970316
$ws.Range("Z13").Value = -3447.4636345431
$ws.Range("AA13").Value = -768.4357733647569
